$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.663.95"
$ws.Range("E2").Value = "  +1.47%  "

$ws.Range("D3").Value = "1.892.09"
$ws.Range("E3").Value = "  +1.98%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'238.58"
$ws.Range("E5").Value = "  +1.28%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.4837"
$ws.Range("E7").Value = "  +1.34%  "

$ws.Range("D8").Value = "'0.2866"
$ws.Range("E8").Value = "  +2.26%  "

$ws.Range("D9").Value = "'0.06558"
$ws.Range("E9").Value = "  +1.48%  "

$ws.Range("D10").Value = "1.911.09"
$ws.Range("E10").Value = "  +2.97%  "

$ws.Range("D11").Value = "'0.07472"
$ws.Range("E11").Value = "  +1.53%  "

$ws.Range("D12").Value = "'16.74"
$ws.Range("E12").Value = "  +3.27%  "

$ws.Range("D13").Value = "'5.101"
$ws.Range("E13").Value = "  +0.26%  "

$ws.Range("D14").Value = "'88.11"
$ws.Range("E14").Value = "  +1.07%  "

$ws.Range("E15").Value = "  +3.40%  "

$ws.Range("D16").Value = "30.630.35"
$ws.Range("E16").Value = "  +1.56%  "

$ws.Range("D17").Value = "'13.26"
$ws.Range("E17").Value = "  +1.12%  "

$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").Value = "2.222.65"
$ws.Range("E19").Value = "  +6.00%  "

$ws.Range("D20").Value = "'0.000007585"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("D21").Value = "'229.94"
$ws.Range("E21").Value = "  +0.74%  "

$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").Value = "'5.268"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").Value = "'6.199"
$ws.Range("E24").Value = "  +1.98%  "

$ws.Range("D25").Value = "'169.74"
$ws.Range("E25").Value = "  +3.82%  "

$ws.Range("D26").Value = "'9.373"
$ws.Range("E26").Value = "  +1.78%  "

$ws.Range("D27").Value = "'18.82"
$ws.Range("E27").Value = "  +1.93%  "

$ws.Range("D28").Value = "'1.960"
$ws.Range("E28").Value = "  +1.80%  "

$ws.Range("D29").Value = "'0.1027"
$ws.Range("E29").Value = "  +11.86%  "

$ws.Range("D30").Value = "'1.399"
$ws.Range("E30").Value = "  -2.67%  "

$ws.Range("D31").Value = "'4.335"
$ws.Range("E31").Value = "  +2.36%  "

$ws.Range("D32").Value = "'4.030"
$ws.Range("E32").Value = "  +1.97%  "

$ws.Range("D33").Value = "'0.05065"
$ws.Range("E33").Value = "  +2.06%  "

$ws.Range("D34").Value = "'1.212"
$ws.Range("E34").Value = "  +5.93%  "

$ws.Range("D35").Value = "'0.7537"
$ws.Range("E35").Value = "  +2.93%  "

$ws.Range("D36").Value = "'1.002"
$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("D37").Value = "'2.716"
$ws.Range("E37").Value = "  +1.11%  "

$ws.Range("D38").Value = "'0.01875"
$ws.Range("E38").Value = "  +1.67%  "

$ws.Range("D40").Value = "'0.9221"
$ws.Range("E40").Value = "  +2.66%  "

$ws.Range("D41").Value = "'2.064"
$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("D42").Value = "'107.15"
$ws.Range("E42").Value = "  +1.04%  "

$ws.Range("D43").Value = "'0.4289"
$ws.Range("E43").Value = "  +1.33%  "

$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("D45").Value = "'5.667"
$ws.Range("E45").Value = "  -4.71%  "

$ws.Range("D46").Value = "'7.424"
$ws.Range("E46").Value = "  +0.83%  "

$ws.Range("D47").Value = "'64.32"
$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("E48").Value = "  -2.73%  "

$ws.Range("D49").Value = "'1.497"
$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("D50").Value = "'8.946"
$ws.Range("E50").Value = "  +2.65%  "

$ws.Range("D51").Value = "'33.98"
$ws.Range("E51").Value = "  +0.33%  "
